# [Admin] (Outlook) Map Signature APIs
# Adds six new rows to the "Snippets" table describing the Outlook
# client-signature APIs (isClientSignatureEnabled[Async], setSignature[Async],
# getComposeType[Async], disableClientSignature[Async]).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New data rows: Class, Member Name, Member ID (methods only), SnippetId, MethodName
$newRows = @(
    @("MessageCompose",     "getComposeType",               2, "outlook-work-with-client-signatures", "getComposeTypeAsync"),
    @("Body",                "setSignatureAsync",            2, "outlook-work-with-client-signatures", "setSignature"),
    @("AppointmentCompose",  "disableClientSignature",       2, "outlook-work-with-client-signatures", "disableClientSignatureAsync"),
    @("MessageCompose",      "disableClientSignature",       2, "outlook-work-with-client-signatures", "disableClientSignatureAsync"),
    @("AppointmentCompose",  "isClientSignatureEnabledAsync",2, "outlook-work-with-client-signatures", "isClientSignatureEnabled"),
    @("MessageCompose",      "isClientSignatureEnabledAsync",2, "outlook-work-with-client-signatures", "isClientSignatureEnabled")
)

foreach ($row in $newRows) {
    $listRow = $lo.ListRows.Add()
    $r = $listRow.Range.Row

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    # Matches the style introduced for the new data (a second cellXfs entry
    # with applyNumberFormat turned on).
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5)).NumberFormat = "General"
}

$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$lastCell = $ws.Cells.Item($lastRow, 1)
$lastCell.Select()
